# CCB May 2023 Push - Process and add 2022 death data to CCB, and update SHA CM 2023
#
# Applies to Sheet1 of the sdohLink workbook:
#   - Rename the "Not_In_School" short-name (J8) to "Not_In_Pre_School" and
#     highlight it in yellow so reviewers notice the rename.
#   - Row 32 ("black"): clear the inSCODA flag (C32) from "x" to a blank
#     space " ".
#   - Add two new SDOH variable rows (33 & 34) for the Index of
#     Concentration at the Extremes (ICE) measures: iceBlack / iceLatino.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 8: enrolled_3_4 -> rename sdohNameShort and flag it yellow -------
$ws.Range("J8").Value = "Not_In_Pre_School"
$ws.Range("J8").Interior.Color = 65535   # RGB(255,255,0) yellow highlight

# --- Row 32: black -> inSCODA flag cleared to a blank space --------------
$ws.Range("C32").Value = " "

# --- Row 33: new ICE_Black variable ---------------------------------------
$ws.Range("A33").Value = "iceBlack [CALCULATED]"
$ws.Range("C33").Value = "x"
$ws.Range("E33").Value = "x"
$ws.Range("H33").Value = "neg"
$ws.Range("I33").Value = "iceBlack"
$ws.Range("J33").Value = "ICE_Black"
$ws.Range("K33").Value = "Index of Concentration at the Extremes (ICE) - Compares population counts of White, NH to Black, NH"

# --- Row 34: new ICE_Latino variable ---------------------------------------
$ws.Range("A34").Value = "iceLatino [CALCULATED]"
$ws.Range("C34").Value = "x"
$ws.Range("E34").Value = "x"
$ws.Range("H34").Value = "neg"
$ws.Range("I34").Value = "iceLatino"
$ws.Range("J34").Value = "ICE_Latino"
$ws.Range("K34").Value = "Index of Concentration at the Extremes (ICE) - Compares population counts of White, NH to Hispanics"

# --- Leave the cursor where the author left it ----------------------------
$ws.Range("K30").Select()
